$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document body:
# "Now, let me just try the local version control use ... -Master"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Now, let me just try the local version control use*") {
        $target = $p
    }
}

$insertPoint = $d.Range($target.Range.End, $target.Range.End)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Build the four new paragraphs (one blank, three with text) as a single
# WordprocessingML fragment so they land after the target paragraph in one
# shot (doing this as separate InsertXML calls causes later inserts to
# collapse into the still-open trailing paragraph).
$newParasXml = (
    "<w:p xmlns:w=`"$wNs`"/>" +
    "<w:p xmlns:w=`"$wNs`"><w:r><w:t>Noz I have a new version again, I need to work on it&gt;</w:t></w:r></w:p>" +
    "<w:p xmlns:w=`"$wNs`"><w:r><w:t>Ther are many modifications from the boss. It will take a long time…</w:t></w:r></w:p>" +
    "<w:p xmlns:w=`"$wNs`"><w:r><w:t>Shit, there are so many hcanges….at one…</w:t></w:r></w:p>"
)

$insertPoint.InsertXML($newParasXml) | Out-Null
